$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Région -> Superficie (km²), Population (habitants)
$data = @(
    @("Auvergne-Rhône-Alpes", 69711, 8030533),
    @("Bourgogne-Franche-Comté", 47784, 2801577),
    @("Bretagne", 27208, 3347004),
    @("Centre-Val de Loire", 39151, 2569510),
    @("Corse", 8680, 342256),
    @("Grand Est", 57433, 5543407),
    @("Guadeloupe", 1703, 383626),
    @("Guyane", 83534, 282107),
    @("Hauts-de-France", 31813, 5995908),
    @("Île-de-France", 12011, 12252917),
    @("Martinique", 1128, 364413),
    @("Mayotte", 2504, 269186),
    @("Normandie", 29906, 3320832),
    @("Nouvelle-Aquitaine", 83809, 5999253),
    @("Occitanie", 72724, 5918981),
    @("Pays de la Loire", 32082, 3800348),
    @("Réunion", 2505, 856547),
    @("Sud-Provence-Alpes-Côte d’Azur", 31400, 5065696)
)

$r = 2
foreach ($row in $data) {
    $bCell = $ws.Cells.Item($r, 2)
    $bCell.NumberFormat = "General"
    $bCell.Value = $row[1]
    $bCell.NumberFormat = "@"

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "General"
    $cCell.Value = $row[2]
    $cCell.NumberFormat = "@"

    $r = $r + 1
}

# Update sheet view: drop the topLeftCell="C1" scroll position and move the
# selection to B2:C19 (the superficie/population data block) instead of D1:Q501.
$ws.Range("B2:C19").Select() | Out-Null
